# Provincial Setup.xlsx update
# Adds the latest data rows to several provincial/territorial sheets,
# updates the Canada sheet's sort range, introduces the new (blue) font
# style used by a couple of rows, and moves the active sheet/selection
# state around to match the saved file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Make sure the new font/style (Calibri 12, Accent1/theme 4 colour) that
# shows up in styles.xml exists. We create it on a throw-away cell far
# outside any sheet's used range and then delete that cell again so the
# workbook ends up with the new <font>/<xf> entries without leaving any
# stray cell behind.
# ---------------------------------------------------------------------
$styleWs = $wb.Worksheets.Item("Canada")
$styleWs.Range("ZZ500").Font.ThemeColor = 5
$null = $styleWs.Range("ZZ500").Delete()

# ---------------------------------------------------------------------
# Atlantic bubble: add row 23
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Atlantic bubble")
$ws.Range("A22:D22").Copy()
$ws.Range("A23:D23").PasteSpecial(-4122)
$ws.Range("A23").Value = 44474
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E25").Select()

# ---------------------------------------------------------------------
# Newfoundland & Labrador: fill in row 22
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Newfoundland & Labrador")
$ws.Range("B21:D21").Copy()
$ws.Range("B22:D22").PasteSpecial(-4122)
$ws.Range("A22").Value = 44454
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("B24").Select()

# ---------------------------------------------------------------------
# New Brunswick: fill rows 16-17 and add row 23
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("New Brunswick")
$ws.Range("A14:D14").Copy()
$ws.Range("A16:D16").PasteSpecial(-4122)
$ws.Range("A16").Value = 44459
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 1

$ws.Range("B16:D16").Copy()
$ws.Range("B17:D17").PasteSpecial(-4122)
$ws.Range("A17").Value = 44474
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 1

$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("F17").Select()

# ---------------------------------------------------------------------
# Canada: add row 20, re-apply the sort, update selection
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Canada")
$ws.Range("A18").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("B19:D19").Copy()
$ws.Range("B20:D20").PasteSpecial(-4122)
$ws.Range("A20").Value = 44446
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1

$sortRange = $ws.Range("A2:D20")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A20"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142
$ws.Sort.Apply()
$ws.Range("E22").Select()

# ---------------------------------------------------------------------
# Alberta: add row 30
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Alberta")
$ws.Range("A29:D29").Copy()
$ws.Range("A30:D30").PasteSpecial(-4122)
$ws.Range("A30").Value = 44455
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = 2
$ws.Range("C30").Select()

# ---------------------------------------------------------------------
# Manitoba: add row 32
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Manitoba")
$ws.Range("A31:D31").Copy()
$ws.Range("A32:D32").PasteSpecial(-4122)
$ws.Range("A32").Value = 44477
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 2
$ws.Range("D32").Value = 1
$ws.Range("G15").Select()

# ---------------------------------------------------------------------
# Saskatchewan: add row 22
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Saskatchewan")
$ws.Range("A21:D21").Copy()
$ws.Range("A22:D22").PasteSpecial(-4122)
$ws.Range("A22").Value = 44456
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 0
$ws.Range("B25").Select()

# ---------------------------------------------------------------------
# Quebec: selection only (full row selected)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Quebec")
$ws.Rows.Item(31).Select()

# ---------------------------------------------------------------------
# British Columbia: add row 21, two-row selection
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("British Columbia")
$ws.Range("A20:D20").Copy()
$ws.Range("A21:D21").PasteSpecial(-4122)
$ws.Range("A21").Value = 44467
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 1
$ws.Range("A20:A21").EntireRow.Select()

# ---------------------------------------------------------------------
# Ontario: add row 27, this is the sheet that ends up active/selected
# so it must be the last sheet we touch.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Ontario")
$ws.Range("A26:D26").Copy()
$ws.Range("A27:D27").PasteSpecial(-4122)
$ws.Range("A27").Value = 44477
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E15").Select()
$ws.Activate()
